# Auto-generated edit script: updates scheduled market-profit values
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7980
$ws.Range("J51").Value = 10142.857
$ws.Range("L51").Value = 10142.857
$ws.Range("N51").Value = -11110.857
$ws.Range("H64").Value = 3368.7058
$ws.Range("I64").Value = 3102
$ws.Range("J64").Value = 3385.375
$ws.Range("K64").Value = 3102
$ws.Range("L64").Value = 3385.375
$ws.Range("M64").Value = -2854
$ws.Range("N64").Value = -3881.375
$ws.Range("H67").Value = 3368.7058
$ws.Range("I67").Value = 3102
$ws.Range("J67").Value = 3385.375
$ws.Range("K67").Value = 3102
$ws.Range("L67").Value = 3385.375
$ws.Range("M67").Value = -2244
$ws.Range("N67").Value = -5101.375
$ws.Range("H137").Value = 1000.29785
$ws.Range("I137").Value = 849.51514
$ws.Range("J137").Value = 1355.7142
$ws.Range("K137").Value = 2548.54542
$ws.Range("L137").Value = 4067.1426
$ws.Range("M137").Value = 1.454580000000078
$ws.Range("N137").Value = -9167.142599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3913.8462
$ws.Range("I63").Value = 2776
$ws.Range("J63").Value = 4625
$ws.Range("K63").Value = 2776
$ws.Range("L63").Value = 4625
$ws.Range("M63").Value = -2090
$ws.Range("N63").Value = -5997
$ws.Range("H66").Value = 3913.8462
$ws.Range("I66").Value = 2776
$ws.Range("J66").Value = 4625
$ws.Range("K66").Value = 13880
$ws.Range("L66").Value = 23125
$ws.Range("M66").Value = -10448
$ws.Range("N66").Value = -29989

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1180327.5
$ws.Range("I134").Value = 1484670.5
$ws.Range("J134").Value = 6433.4287
$ws.Range("K134").Value = 4454011.5
$ws.Range("L134").Value = 19300.2861
$ws.Range("M134").Value = -4451476.5
$ws.Range("N134").Value = -24370.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 864.2857
$ws.Range("I16").Value = 712.5
$ws.Range("J16").Value = 1066.6666
$ws.Range("K16").Value = 712.5
$ws.Range("L16").Value = 1066.6666
$ws.Range("M16").Value = -425.5
$ws.Range("N16").Value = -1640.6666
$ws.Range("H86").Value = 18286.715
$ws.Range("I86").Value = 6335.6665
$ws.Range("J86").Value = 27250
$ws.Range("K86").Value = 6335.6665
$ws.Range("L86").Value = 27250
$ws.Range("M86").Value = -5212.6665
$ws.Range("N86").Value = -29496
$ws.Range("H89").Value = 18286.715
$ws.Range("I89").Value = 6335.6665
$ws.Range("J89").Value = 27250
$ws.Range("K89").Value = 31678.3325
$ws.Range("L89").Value = 136250
$ws.Range("M89").Value = -26062.3325
$ws.Range("N89").Value = -147482
$ws.Range("H105").Value = 1300
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 1300
$ws.Range("N105").Value = -4794
$ws.Range("M105").Value = ""
$ws.Range("H113").Value = 864.2857
$ws.Range("I113").Value = 712.5
$ws.Range("J113").Value = 1066.6666
$ws.Range("K113").Value = 712.5
$ws.Range("L113").Value = 1066.6666
$ws.Range("M113").Value = 1457.5
$ws.Range("N113").Value = -5406.6666
$ws.Range("H118").Value = 31400
$ws.Range("J118").Value = 31400
$ws.Range("L118").Value = 31400
$ws.Range("N118").Value = -34714
$ws.Range("H119").Value = 31333.334
$ws.Range("J119").Value = 31333.334
$ws.Range("L119").Value = 31333.334
$ws.Range("N119").Value = -41009.334
$ws.Range("H141").Value = 61828.75
$ws.Range("J141").Value = 61828.75
$ws.Range("L141").Value = 61828.75
$ws.Range("N141").Value = -72188.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 853.2778
$ws.Range("I5").Value = 384.2143
$ws.Range("J5").Value = 2495
$ws.Range("K5").Value = 1152.6429
$ws.Range("L5").Value = 7485
$ws.Range("M5").Value = -1040.6429
$ws.Range("N5").Value = -7709
$ws.Range("H39").Value = 1990.5555
$ws.Range("J39").Value = 5200
$ws.Range("L39").Value = 15600
$ws.Range("N39").Value = -16188
$ws.Range("H51").Value = 3844.8147
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 3877.3076
$ws.Range("K51").Value = 9000
$ws.Range("L51").Value = 11631.9228
$ws.Range("M51").Value = -8540
$ws.Range("N51").Value = -12551.9228
$ws.Range("H55").Value = 754.375
$ws.Range("I55").Value = 433.33334
$ws.Range("J55").Value = 828.46155
$ws.Range("K55").Value = 1300.00002
$ws.Range("L55").Value = 2485.38465
$ws.Range("M55").Value = -1123.00002
$ws.Range("N55").Value = -2839.38465
$ws.Range("H122").Value = 45259.68
$ws.Range("J122").Value = 2087.818
$ws.Range("L122").Value = 18790.362
$ws.Range("N122").Value = -23690.362
$ws.Range("H131").Value = 908.14
$ws.Range("I131").Value = 307.6
$ws.Range("J131").Value = 939.7474
$ws.Range("K131").Value = 922.8000000000001
$ws.Range("L131").Value = 2819.2422
$ws.Range("M131").Value = 4117.2
$ws.Range("N131").Value = -12899.2422
$ws.Range("H135").Value = 853.2778
$ws.Range("I135").Value = 384.2143
$ws.Range("J135").Value = 2495
$ws.Range("K135").Value = 3457.9287
$ws.Range("L135").Value = 22455
$ws.Range("M135").Value = -922.9286999999999
$ws.Range("N135").Value = -27525

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2153.8462
$ws.Range("I80").Value = 2133.3333
$ws.Range("J80").Value = 2400
$ws.Range("K80").Value = 2133.3333
$ws.Range("L80").Value = 2400
$ws.Range("M80").Value = -1135.3333
$ws.Range("N80").Value = -4396
$ws.Range("H83").Value = 2153.8462
$ws.Range("I83").Value = 2133.3333
$ws.Range("J83").Value = 2400
$ws.Range("K83").Value = 10666.6665
$ws.Range("L83").Value = 12000
$ws.Range("M83").Value = -5674.666499999999
$ws.Range("N83").Value = -21984
$ws.Range("H113").Value = 1694.0667
$ws.Range("I113").Value = 988.875
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 988.875
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 1181.125
$ws.Range("N113").Value = -6840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 624.95654
$ws.Range("I55").Value = 287.8889
$ws.Range("J55").Value = 841.6429000000001
$ws.Range("K55").Value = 287.8889
$ws.Range("L55").Value = 841.6429000000001
$ws.Range("M55").Value = -114.8889
$ws.Range("N55").Value = -1187.6429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1629.8334
$ws.Range("J81").Value = 2400
$ws.Range("L81").Value = 4800
$ws.Range("N81").Value = -6922
$ws.Range("H84").Value = 1629.8334
$ws.Range("J84").Value = 2400
$ws.Range("L84").Value = 24000
$ws.Range("N84").Value = -34608
$ws.Range("H113").Value = 930.5
$ws.Range("I113").Value = 725
$ws.Range("J113").Value = 1067.5
$ws.Range("K113").Value = 2175
$ws.Range("L113").Value = 3202.5
$ws.Range("M113").Value = -5
$ws.Range("N113").Value = -7542.5
